# Update column B (word) values for rows 2-193 with the new working set of sequences
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    'nutzen',
    'schwingen',
    'doppeln',
    'schnellen',
    'hacken',
    'schwächen',
    'gründen',
    'boxen',
    'buchen',
    'wirken',
    'erben',
    'fällen',
    'kennen',
    'klingeln',
    'fangen',
    'bellen',
    'lesen',
    'pflanzen',
    'äußern',
    'tropfen',
    'schauen',
    'treiben',
    'laufen',
    'trennen',
    'dringen',
    'spüren',
    'münzen',
    'heben',
    'zählen',
    'regeln',
    'stehlen',
    'wachsen',
    'gelten',
    'dienen',
    'rasen',
    'holen',
    'regnen',
    'schneiden',
    'fallen',
    'platzen',
    'graben',
    'loben',
    'ändern',
    'kehren',
    'krachen',
    'kichern',
    'seufzen',
    'altern',
    'wärmen',
    'foltern',
    'betteln',
    'lehnen',
    'führen',
    'bauen',
    'flüstern',
    'heilen',
    'wehtun',
    'weichen',
    'kosten',
    'helfen',
    'knarren',
    'strahlen',
    'schlucken',
    'schmecken',
    'fahren',
    'bremsen',
    'wenden',
    'landen',
    'wehren',
    'filmen',
    'pressen',
    'hauen',
    'scheinen',
    'drehen',
    'trotzen',
    'töten',
    'klettern',
    'stören',
    'schwören',
    'runden',
    'fischen',
    'pfeifen',
    'feiern',
    'werden',
    'drohen',
    'schreiten',
    'sperren',
    'brauchen',
    'werfen',
    'liegen',
    'wahren',
    'planen',
    'folgen',
    'machen',
    'kürzen',
    'saufen',
    'fließen',
    'opfern',
    'heulen',
    'schulden',
    'treffen',
    'ärgern',
    'sichern',
    'reiten',
    'siegen',
    'stillen',
    'streichen',
    'sitzen',
    'deuten',
    'knurren',
    'tollen',
    'gnaden',
    'wundern',
    'jagen',
    'sprengen',
    'sinken',
    'jubeln',
    'suchen',
    'meistern',
    'grüßen',
    'rufen',
    'hoffen',
    'stecken',
    'pissen',
    'melden',
    'zielen',
    'backen',
    'irren',
    'reizen',
    'mauern',
    'leihen',
    'freuen',
    'bergen',
    'klingen',
    'ruhen',
    'quälen',
    'fügen',
    'flüchten',
    'morden',
    'geben',
    'lügen',
    'sprechen',
    'binden',
    'trauen',
    'hören',
    'scheitern',
    'scheiden',
    'achten',
    'kümmern',
    'flehen',
    'sorgen',
    'bluten',
    'nerven',
    'meinen',
    'liefern',
    'nähen',
    'zeichnen',
    'schenken',
    'wandern',
    'stürzen',
    'decken',
    'parken',
    'zünden',
    'biegen',
    'warnen',
    'lockern',
    'stürmen',
    'leiden',
    'zögern',
    'enden',
    'arten',
    'räumen',
    'stammen',
    'schwimmen',
    'zeigen',
    'bilden',
    'spinnen',
    'zeugen',
    'schrecken',
    'greifen',
    'malen',
    'wüten',
    'schützen',
    'bitten',
    'orten',
    'spielen',
    'wecken',
    'garen',
    'sterben',
    'mögen',
    'formen',
    'ehren'
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $words[$i]
}

